$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry two extra (empty) header cells in C1:D1 that only
# held formatting. Drop them entirely so the used range shrinks back down
# to A1:B2, matching the trimmed "level" import/export template.
$ws.Range("C1:D1").Clear()

# Re-assert the bold header formatting on the two remaining header cells.
$ws.Range("A1:B1").Font.Bold = $true

# Row 2 becomes the sample/seed data row for the Level template:
# level_kode = "KRU", level_nama = "Kurir" (was the numeric placeholder 5).
$ws.Range("A2").Value = "KRU"
$ws.Range("B2").Value = "Kurir"

# Leave the cursor where the author left it after editing.
$ws.Range("D14").Select() | Out-Null
